$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 5849.9
$ws.Range("J70").Value = 6687.5
$ws.Range("L70").Value = 20062.5
$ws.Range("N70").Value = -20602.5
$ws.Range("H73").Value = 5849.9
$ws.Range("J73").Value = 6687.5
$ws.Range("L73").Value = 20062.5
$ws.Range("N73").Value = -21934.5
$ws.Range("H80").Value = 1590
$ws.Range("I80").Value = 699.3333
$ws.Range("J80").Value = 2258
$ws.Range("K80").Value = 2097.9999
$ws.Range("L80").Value = 6774
$ws.Range("M80").Value = -1099.9999
$ws.Range("N80").Value = -8770
$ws.Range("H83").Value = 1590
$ws.Range("I83").Value = 699.3333
$ws.Range("J83").Value = 2258
$ws.Range("K83").Value = 6293.9997
$ws.Range("L83").Value = 20322
$ws.Range("M83").Value = -1301.9997
$ws.Range("N83").Value = -30306
$ws.Range("H96").Value = 1634.2778
$ws.Range("I96").Value = 417.86667
$ws.Range("K96").Value = 1253.60001
$ws.Range("M96").Value = 119.3999899999999
$ws.Range("H99").Value = 6678.467
$ws.Range("J99").Value = 12880.857
$ws.Range("L99").Value = 38642.571
$ws.Range("N99").Value = -41638.571
$ws.Range("H101").Value = 3580.5
$ws.Range("I101").Value = 551
$ws.Range("J101").Value = 5398.2
$ws.Range("K101").Value = 1653
$ws.Range("L101").Value = 16194.6
$ws.Range("M101").Value = -31
$ws.Range("N101").Value = -19438.6
$ws.Range("H113").Value = 3526.7
$ws.Range("J113").Value = 3879
$ws.Range("L113").Value = 3879
$ws.Range("N113").Value = -10387
$ws.Range("H132").Value = 9490.528
$ws.Range("I132").Value = 1660.4186
$ws.Range("K132").Value = 4981.2558
$ws.Range("M132").Value = -2451.2558
$ws.Range("H137").Value = 12719550
$ws.Range("I137").Value = 834794.4399999999
$ws.Range("K137").Value = 2504383.32
$ws.Range("M137").Value = -2501833.32
$ws.Range("H138").Value = 3589.798
$ws.Range("I138").Value = 1166.1724
$ws.Range("K138").Value = 3498.5172
$ws.Range("M138").Value = 1641.4828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14456.266
$ws.Range("I32").Value = 14256.096
$ws.Range("J32").Value = 15323.667
$ws.Range("K32").Value = 14256.096
$ws.Range("L32").Value = 15323.667
$ws.Range("M32").Value = -13969.096
$ws.Range("N32").Value = -15897.667
$ws.Range("H45").Value = 3366
$ws.Range("J45").Value = 4499
$ws.Range("L45").Value = 4499
$ws.Range("N45").Value = -5253
$ws.Range("H61").Value = 14257.571
$ws.Range("I61").Value = 19345
$ws.Range("J61").Value = 5100.2
$ws.Range("K61").Value = 19345
$ws.Range("L61").Value = 5100.2
$ws.Range("M61").Value = -19133
$ws.Range("N61").Value = -5524.2
$ws.Range("H74").Value = 1560.4
$ws.Range("I74").Value = 1233.0714
$ws.Range("J74").Value = 1977
$ws.Range("K74").Value = 1233.0714
$ws.Range("L74").Value = 1977
$ws.Range("M74").Value = -359.0714
$ws.Range("N74").Value = -3725
$ws.Range("H77").Value = 1560.4
$ws.Range("I77").Value = 1233.0714
$ws.Range("J77").Value = 1977
$ws.Range("K77").Value = 6165.357
$ws.Range("L77").Value = 9885
$ws.Range("M77").Value = -1797.357
$ws.Range("N77").Value = -18621
$ws.Range("H102").Value = 508954.62
$ws.Range("I102").Value = 549351.4
$ws.Range("K102").Value = 549351.4
$ws.Range("M102").Value = -547729.4
$ws.Range("H132").Value = 14595.28
$ws.Range("I132").Value = 27785.092
$ws.Range("J132").Value = 4231.857
$ws.Range("K132").Value = 83355.276
$ws.Range("L132").Value = 12695.571
$ws.Range("M132").Value = -80825.276
$ws.Range("N132").Value = -17755.571
$ws.Range("H136").Value = 14257.571
$ws.Range("I136").Value = 19345
$ws.Range("J136").Value = 5100.2
$ws.Range("K136").Value = 58035
$ws.Range("L136").Value = 15300.6
$ws.Range("M136").Value = -55485
$ws.Range("N136").Value = -20400.6
$ws.Range("H138").Value = 63143
$ws.Range("J138").Value = 63143
$ws.Range("L138").Value = 63143
$ws.Range("N138").Value = -73423

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2084952.2
$ws.Range("I99").Value = 2605690.2
$ws.Range("K99").Value = 2605690.2
$ws.Range("M99").Value = -2604192.2
$ws.Range("H105").Value = 2748.75
$ws.Range("H107").Value = 969
$ws.Range("I107").Value = 969
$ws.Range("K107").Value = 969
$ws.Range("M107").Value = 951

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6455.976
$ws.Range("J31").Value = 7747.9395
$ws.Range("L31").Value = 7747.9395
$ws.Range("N31").Value = -8337.9395
$ws.Range("H34").Value = 6455.976
$ws.Range("J34").Value = 7747.9395
$ws.Range("L34").Value = 7747.9395
$ws.Range("N34").Value = -8151.9395
$ws.Range("H58").Value = 316386.4
$ws.Range("I58").Value = 1854.8334
$ws.Range("J58").Value = 505105.34
$ws.Range("K58").Value = 1854.8334
$ws.Range("L58").Value = 505105.34
$ws.Range("M58").Value = -1651.8334
$ws.Range("N58").Value = -505511.34
$ws.Range("H99").Value = 7363.5654
$ws.Range("I99").Value = 5241.7144
$ws.Range("K99").Value = 5241.7144
$ws.Range("M99").Value = -3743.7144
$ws.Range("H122").Value = 3216.516
$ws.Range("I122").Value = 1685.15
$ws.Range("J122").Value = 6000.8184
$ws.Range("K122").Value = 5055.450000000001
$ws.Range("L122").Value = 18002.4552
$ws.Range("M122").Value = -2605.450000000001
$ws.Range("N122").Value = -22902.4552
$ws.Range("H126").Value = 7363.5654
$ws.Range("I126").Value = 5241.7144
$ws.Range("K126").Value = 15725.1432
$ws.Range("M126").Value = -13255.1432
$ws.Range("H132").Value = 10762873
$ws.Range("I132").Value = 13901342
$ws.Range("J132").Value = 2407.2856
$ws.Range("K132").Value = 41704026
$ws.Range("L132").Value = 7221.8568
$ws.Range("M132").Value = -41701496
$ws.Range("N132").Value = -12281.8568
$ws.Range("H134").Value = 1470.7894
$ws.Range("I134").Value = 1546.4412
$ws.Range("K134").Value = 4639.3236
$ws.Range("M134").Value = -2104.3236
$ws.Range("H136").Value = 316386.4
$ws.Range("I136").Value = 1854.8334
$ws.Range("J136").Value = 505105.34
$ws.Range("K136").Value = 5564.5002
$ws.Range("L136").Value = 1515316.02
$ws.Range("M136").Value = -3014.5002
$ws.Range("N136").Value = -1520416.02
$ws.Range("H141").Value = 85824.12
$ws.Range("J141").Value = 87419.625
$ws.Range("L141").Value = 87419.625
$ws.Range("N141").Value = -97779.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3729646
$ws.Range("I4").Value = 706306.5
$ws.Range("J4").Value = 14009000
$ws.Range("K4").Value = 2118919.5
$ws.Range("L4").Value = 42027000
$ws.Range("M4").Value = -2118807.5
$ws.Range("N4").Value = -42027224
$ws.Range("H8").Value = 73.875
$ws.Range("I8").Value = 73.875
$ws.Range("K8").Value = 221.625
$ws.Range("M8").Value = -82.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 755.3158
$ws.Range("I97").Value = 596.8
$ws.Range("J97").Value = 1349.75
$ws.Range("K97").Value = 596.8
$ws.Range("L97").Value = 1349.75
$ws.Range("M97").Value = -100.8
$ws.Range("N97").Value = -2341.75
$ws.Range("H102").Value = 23818942
$ws.Range("I102").Value = 38471596
$ws.Range("K102").Value = 38471596
$ws.Range("M102").Value = -38469974
$ws.Range("H122").Value = 921643.4399999999
$ws.Range("I122").Value = 1835625
$ws.Range("J122").Value = 7661.8335
$ws.Range("K122").Value = 5506875
$ws.Range("L122").Value = 22985.5005
$ws.Range("M122").Value = -5504425
$ws.Range("N122").Value = -27885.5005
$ws.Range("H123").Value = 53753
$ws.Range("J123").Value = 53753
$ws.Range("L123").Value = 53753
$ws.Range("N123").Value = -58653
$ws.Range("H126").Value = 4332.25
$ws.Range("I126").Value = 1961.3077
$ws.Range("K126").Value = 5883.9231
$ws.Range("M126").Value = -3413.9231
$ws.Range("H132").Value = 604986.5600000001
$ws.Range("I132").Value = 169978.08
$ws.Range("J132").Value = 1257499.2
$ws.Range("K132").Value = 509934.24
$ws.Range("L132").Value = 3772497.6
$ws.Range("M132").Value = -507404.24
$ws.Range("N132").Value = -3777557.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 79898.5
$ws.Range("J6").Value = 79898.5
$ws.Range("L6").Value = 79898.5
$ws.Range("N6").Value = -80122.5
$ws.Range("H87").Value = 42000
$ws.Range("J87").Value = 42000
$ws.Range("L87").Value = 42000
$ws.Range("N87").Value = -44246
$ws.Range("H90").Value = 42000
$ws.Range("J90").Value = 42000
$ws.Range("L90").Value = 126000
$ws.Range("N90").Value = -137232
$ws.Range("H93").Value = 2055.5625
$ws.Range("I93").Value = 2142.0715
$ws.Range("J93").Value = 1450
$ws.Range("K93").Value = 2142.0715
$ws.Range("L93").Value = 1450
$ws.Range("M93").Value = -894.0715
$ws.Range("N93").Value = -3946
$ws.Range("H132").Value = 6900.091
$ws.Range("J132").Value = 4934.75
$ws.Range("L132").Value = 14804.25
$ws.Range("N132").Value = -19864.25
$ws.Range("H136").Value = 3323.36
$ws.Range("I136").Value = 2188.6843
$ws.Range("K136").Value = 6566.0529
$ws.Range("M136").Value = -4016.0529
$ws.Range("H140").Value = 67181
$ws.Range("J140").Value = 67181
$ws.Range("L140").Value = 67181
$ws.Range("N140").Value = -77541

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3562.9473
$ws.Range("I122").Value = 3562.9473
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10688.8419
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8238.841899999999
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 20001300
$ws.Range("I132").Value = 1591.0769
$ws.Range("J132").Value = 41667650
$ws.Range("K132").Value = 4773.2307
$ws.Range("L132").Value = 125002950
$ws.Range("M132").Value = -2243.2307
$ws.Range("N132").Value = -125008010
